## Replace the hard-coded "www.drpaulduenas.com" footer text with a
## configurable MERGEFIELD ("=website") field, matching the structure
## Word produces for Insert > Quick Parts > Field > MergeField:
##
##   <w:r><w:fldChar w:fldCharType="begin"/></w:r>
##   <w:r><w:instrText> MERGEFIELD =website \* MERGEFORMAT </w:instrText></w:r>
##   <w:r><w:fldChar w:fldCharType="separate"/></w:r>
##   <w:r><w:t>&#171;=website&#187;</w:t></w:r>
##   <w:r><w:fldChar w:fldCharType="end"/></w:r>
##
## all sharing the run formatting (Avenir Book, bold, 20 half-points)
## the literal text run used to carry.

$d = $word.ActiveDocument

# The text lives in the document's default (primary) footer.
$footer = $d.Sections.Item(1).Footers.Item(1)

$target = "www.drpaulduenas.com"
$rng = $footer.Range.Duplicate
$found = $rng.Find.Execute($target, $false, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find '$target' in the default footer"
}

# Drop the literal-text run; InsertXML below will repopulate the
# (now collapsed) range with the field-code run sequence.
$rng.Delete()

$guillemetOpen = [char]0x00AB
$guillemetClose = [char]0x00BB

$fieldXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
  'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' +
  'w14:paraId="3CBC13F5" w14:textId="77777777" w:rsidR="00B55E9A" ' +
  'w:rsidRDefault="00B55E9A" w:rsidP="00B55E9A">' +
    '<w:pPr>' +
      '<w:pStyle w:val="Footer"/>' +
      '<w:jc w:val="center"/>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/>' +
        '<w:b/>' +
        '<w:sz w:val="20"/>' +
        '<w:szCs w:val="20"/>' +
      '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/>' +
        '<w:b/>' +
        '<w:sz w:val="20"/>' +
        '<w:szCs w:val="20"/>' +
      '</w:rPr>' +
      '<w:fldChar w:fldCharType="begin"/>' +
    '</w:r>' +
    '<w:r>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/>' +
        '<w:b/>' +
        '<w:sz w:val="20"/>' +
        '<w:szCs w:val="20"/>' +
      '</w:rPr>' +
      '<w:instrText xml:space="preserve"> MERGEFIELD =website \* MERGEFORMAT </w:instrText>' +
    '</w:r>' +
    '<w:r>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/>' +
        '<w:b/>' +
        '<w:sz w:val="20"/>' +
        '<w:szCs w:val="20"/>' +
      '</w:rPr>' +
      '<w:fldChar w:fldCharType="separate"/>' +
    '</w:r>' +
    '<w:r>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/>' +
        '<w:b/>' +
        '<w:sz w:val="20"/>' +
        '<w:szCs w:val="20"/>' +
      '</w:rPr>' +
      "<w:t>$guillemetOpen=website$guillemetClose</w:t>" +
    '</w:r>' +
    '<w:r>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/>' +
        '<w:b/>' +
        '<w:sz w:val="20"/>' +
        '<w:szCs w:val="20"/>' +
      '</w:rPr>' +
      '<w:fldChar w:fldCharType="end"/>' +
    '</w:r>' +
  '</w:p>'

$rng.InsertXML($fieldXml)

Write-Output "Replaced '$target' with MERGEFIELD =website in the default footer."
